# Generate Report for Handback
# Updates the handback status (now out of sync with en-US) and refreshes the
# de-de "Correspond Handback DateTime" for the 56da58e0... file after a new
# handback round was generated.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: not in sync with en-US"

# --- Overview sheet: zh-cn / de-de status columns for both files ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Columns E & F grow to fit the new (longer) status text
$wsOverview.Columns.Item(5).ColumnWidth = 32.65
$wsOverview.Columns.Item(6).ColumnWidth = 32.65

# --- zh-cn sheet: Status column ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Columns.Item(3).ColumnWidth = 32.65

# --- de-de sheet: Status column + refreshed handback datetime ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Columns.Item(3).ColumnWidth = 32.65

# New handback was generated for 56da58e0...de-de.xlf at 17:41:43
$wsDeDe.Range("K2").Value = "2016-10-19 17:41:43"
